# Generate Report for Handoff
#
# For the most recently added localized file in each language sheet
# (row 5 = the "b01c9cac-..." file), stamp the "Latest Handoff Datetime"
# column (D) with the datetime the report was generated, per language.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-03 07:19:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-03 07:19:11"
